# Generate Report for Handoff
#
# A new handoff happened for "b" (source b.md / blob sha
# 63290e5768f688058c7b37413b0a5c26c308f864) in both the zh-cn and de-de
# locales. Update the per-locale status rows and the Overview roll-up row
# to reflect the new "Ready for handoff" status and the new handoff
# file/datetime values.

$wb = $excel.ActiveWorkbook

$statusText    = "Ready for handoff"
$overviewDate  = "2016-29-13 00:29:21"

$zhFile        = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhDatetime    = "2016-03-13 00:29:17"

$deFile        = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$deDatetime    = "2016-03-13 00:29:21"

# --- Overview sheet: row 3 is the "b.md" row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText
$wsOverview.Range("D3").Value = $overviewDate

# --- zh-cn sheet: row 3 is the "b.md" row ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusText
$wsZhCn.Range("D3").Value = $zhFile
$wsZhCn.Range("E3").Value = $zhDatetime

# --- de-de sheet: row 3 is the "b.md" row ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusText
$wsDeDe.Range("D3").Value = $deFile
$wsDeDe.Range("E3").Value = $deDatetime
